# Module 2 update: emphasize the "01 | C# features add productivity and
# conciseness" row of the agenda table on slide 6 by making its text bold
# and italic.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)
$tbl = $shp.Table

$cell = $tbl.Cell(2, 1)
$tr = $cell.Shape.TextFrame.TextRange
$tr.Font.Bold = $true
$tr.Font.Italic = $true
